# "Valid Species" taxonomy clean-up:
#  - add "Pogoniulus atroflavus" as a new valid species, in alphabetical
#    order just before "Pogoniulus bilineatus" (row 3160)
#  - remove the "Turdus philomelos" (Song Thrush) group of 5 entries
#    (species + 4 subspecies), originally at rows 4136-4140
#  - fix the duplicated "Sterna bergii bergii" entry (originally row 3724)
#    to read "Sterna bergii cristata"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valid Species")

# 1) Insert the new species row just above "Pogoniulus bilineatus" (row 3160).
$ws.Range("A3160").EntireRow.Insert() | Out-Null
$ws.Range("A3160").Value = "Pogoniulus atroflavus"

# 2) Remove the "Turdus philomelos" group. The insert above shifted these
#    rows down by one: 4136-4140 -> 4137-4141.
$ws.Range("A4137:A4141").EntireRow.Delete() | Out-Null

# 3) Fix the duplicated row. It was row 3724 originally; the insert in step
#    1 (which happened above/before it) shifted it down to row 3725, and
#    the later deletion (which happened well below it) doesn't affect it.
$ws.Range("A3725").Value = "Sterna bergii cristata"
